$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.066.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -5.92%  '

$ws.Range("D3").Value = "'2.445.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -8.71%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = "'538.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.82%  '

$ws.Range("D6").Value = "'145.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.44%  '

$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.26%  '

$ws.Range("E8").Value = '  -2.64%  '

$ws.Range("D9").Value = "'2.457.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.33%  '

$ws.Range("D10").Value = "'0.0991"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.02%  '

$ws.Range("E11").Value = '  -1.38%  '

$ws.Range("D12").Value = "'5.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.67%  '

$ws.Range("E13").Value = '  -4.25%  '

$ws.Range("D14").Value = "'2.885.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.50%  '

$ws.Range("D15").Value = "'23.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.71%  '

$ws.Range("D16").Value = "'58.958.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.96%  '

$ws.Range("D17").Value = "'0.0000138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.53%  '

$ws.Range("D18").Value = "'2.499.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.70%  '

$ws.Range("E19").Value = '  -5.68%  '

$ws.Range("E20").Value = '  -5.25%  '

$ws.Range("D21").Value = "'323.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.06%  '

$ws.Range("E22").Value = '  -3.36%  '

$ws.Range("D23").Value = "'5.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.69%  '

$ws.Range("D24").Value = "'60.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.96%  '

$ws.Range("D25").Value = "'0.449"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -11.83%  '

$ws.Range("E26").Value = '  -5.07%  '

$ws.Range("E27").Value = '  -2.39%  '

$ws.Range("D28").Value = "'7.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.01%  '

$ws.Range("D29").Value = "'6.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.61%  '

$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Value = "'1.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.99%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'1.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.61%  '

$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = "'0.0₃0771"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.95%  '

$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.11%  '

$ws.Range("D34").Value = "'157.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.31%  '

$ws.Range("D35").Value = "'1.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.90%  '

$ws.Range("D36").Value = "'18.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.21%  '

$ws.Range("E37").Value = '  -8.37%  '

$ws.Range("E38").Value = '  -4.16%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = "'313.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.60%  '

$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = "'5.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.99%  '

$ws.Range("E41").Value = '  -5.28%  '

$ws.Range("D42").Value = "'0.834"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.94%  '

$ws.Range("E43").Value = '  -6.75%  '

$ws.Range("D44").Value = "'0.995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.30%  '

$ws.Range("E45").Value = '  -2.44%  '

$ws.Range("E46").Value = '  -4.02%  '

$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").Value = "'0.0526"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.71%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = "'0.0933"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.53%  '

$ws.Range("D49").Value = "'123.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.95%  '

$ws.Range("E50").Value = '  -4.54%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = "'18.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.51%  '
